$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped by
# one day (46074 -> 46075) for every data row (rows 2 through 372).
$startRow = 2
$endRow = 372
$newValue = 46075

for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
